$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New measurement rows to append (distance, consume, speed, temp_inside, temp_outside, specials, gas_type, AC, rain)
$rows = @(
    @(130.30000000000001, 4.5999999999999996, 85, 22, 12, $null,                     "E10",  0, 0),
    @(67.2,               4.3,                67, 22, 18, $null,                     "E10",  0, 0),
    @(43.7,               4.7,                44, 22, 9,  "half rain half sun",      "SP98", 0, 1),
    @(12.1,               4.2,                43, 22, 4,  $null,                     "SP98", 0, 0),
    @(56.1,               4.8,                82, 22, 13, $null,                     "SP98", 0, 0),
    @(39,                 4.0999999999999996, 61, 22, 16, $null,                     "SP98", 0, 0)
)

$startRow = 191
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    if ($data[5] -ne $null) {
        $ws.Cells.Item($r, 6).Value = $data[5]
    }
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

# Resize the table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I196"))

# Update selection to mirror the authored workbook state
$ws.Range("D197").Select()
